# Ignore columns whose header starts with "!" -- demo/test data added to
# column G of the "First" sheet, plus a header rename on A1.
#
# Shared strings are appended in the order they are first encountered, so
# we populate column G (top to bottom) before touching A1 to keep the
# resulting xl/sharedStrings.xml ordering stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First")

$ws.Range("G1").Value  = "!ddasda"
$ws.Range("G3").Value  = "fsd"
$ws.Range("G5").Value  = "fdf"
$ws.Range("G8").Value  = "fsd"
$ws.Range("G9").Value  = "fsd"
$ws.Range("G11").Value = "fsdf"
$ws.Range("G13").Value = "ddd"
$ws.Range("G14").Value = "ddd"

$ws.Range("A1").Value = "!a"
